$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move Thiago (previously row 6) up to row 2, columns A:B
$ws.Range("A2").Value = "Thiago"
$ws.Range("B2").Value = "thiago.honorato.pb@compasso.com.br"

# Move Amanda/Bruna/Luciana/Nathalia to columns H:I (rows 2-5),
# with the name cells getting an explicit black font color.
$ws.Range("H2").Value = "Amanda"
$ws.Range("I2").Value = "amanda.wanderley.pb@compasso.com.br"
$ws.Range("H2").Font.Color = 0

$ws.Range("H3").Value = "Bruna"
$ws.Range("I3").Value = "Maria.Nunes.pb@compasso.com.br"
$ws.Range("H3").Font.Color = 0

$ws.Range("H4").Value = "Luciana"
$ws.Range("I4").Value = "luciana.maciel.pb@compasso.com.br"
$ws.Range("H4").Font.Color = 0

$ws.Range("H5").Value = "Nathalia"
$ws.Range("I5").Value = "bruna.souza.pb@compasso.com.br"
$ws.Range("H5").Font.Color = 0

# Clear the old A3:B6 block which is no longer used.
$ws.Range("A3:B6").Clear()

# Update selection to match the saved state.
$ws.Range("B10").Select()
